{"js": "// Auto-generated: replace each unique old value with its new value.\n// Each (old, new) pair appears exactly once in the document body, so a targeted\n// body.search(...) + insertText(..., 'Replace') pair is safe and unambiguous.\nconst replacements = [\n  [\"2024-07-10 Wednesday\", \"2024-07-11 Thursday\"],\n  [\"33\u00f78=4, 1\", \"39\u00f72=19, 1\"],\n  [\"31\u00f76=5, 1\", \"30\u00f73=10, 0\"],\n  [\"94\u00f77=13, 3\", \"85\u00f73=28, 1\"],\n  [\"41\u00f78=5, 1\", \"93\u00f79=10, 3\"],\n  [\"22\u00f72=11, 0\", \"67\u00f78=8, 3\"],\n  [\"94\u00f76=15, 4\", \"82\u00f73=27, 1\"],\n  [\"97\u00f72=48, 1\", \"65\u00f75=13, 0\"],\n  [\"33\u00f76=5, 3\", \"83\u00f77=11, 6\"],\n  [\"65\u00f79=7, 2\", \"26\u00f79=2, 8\"],\n  [\"60\u00f73=20, 0\", \"57\u00f73=19, 0\"],\n  [\"58\u00f76=9, 4\", \"73\u00f74=18, 1\"],\n  [\"86\u00f73=28, 2\", \"46\u00f74=11, 2\"],\n  [\"53\u00f79=5, 8\", \"12\u00f75=2, 2\"],\n  [\"95\u00f73=31, 2\", \"95\u00f72=47, 1\"],\n  [\"67\u00f72=33, 1\", \"61\u00f73=20, 1\"],\n  [\"63\u00f79=7, 0\", \"17\u00f79=1, 8\"],\n  [\"98\u00f79=10, 8\", \"86\u00f79=9, 5\"],\n  [\"81\u00f79=9, 0\", \"24\u00f77=3, 3\"],\n  [\"58\u00f73=19, 1\", \"75\u00f74=18, 3\"],\n  [\"23\u00f79=2, 5\", \"56\u00f75=11, 1\"],\n  [\"93\u00f74=23, 1\", \"23\u00f79=2, 5\"],\n  [\"53\u00f72=26, 1\", \"93\u00f78=11, 5\"],\n  [\"28\u00f75=5, 3\", \"10\u00f73=3, 1\"],\n  [\"74\u00f76=12, 2\", \"27\u00f72=13, 1\"],\n  [\"72\u00f78=9, 0\", \"49\u00f76=8, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for ${JSON.stringify(oldText)}, found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, 'Replace');\n}\n\nawait context.sync();\n", "ps1": "# Auto-generated: replace each unique old value with its new value using\n# Word's native Find/Replace (Range.Find.Execute), mirroring the diff 1:1.\n# Each (old, new) pair appears exactly once in the document, so Replace:=wdReplaceAll\n# (2) only ever touches the single intended occurrence.\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @(\"2024-07-10 Wednesday\", \"2024-07-11 Thursday\"),\n    @(\"33\u00f78=4, 1\", \"39\u00f72=19, 1\"),\n    @(\"31\u00f76=5, 1\", \"30\u00f73=10, 0\"),\n    @(\"94\u00f77=13, 3\", \"85\u00f73=28, 1\"),\n    @(\"41\u00f78=5, 1\", \"93\u00f79=10, 3\"),\n    @(\"22\u00f72=11, 0\", \"67\u00f78=8, 3\"),\n    @(\"94\u00f76=15, 4\", \"82\u00f73=27, 1\"),\n    @(\"97\u00f72=48, 1\", \"65\u00f75=13, 0\"),\n    @(\"33\u00f76=5, 3\", \"83\u00f77=11, 6\"),\n    @(\"65\u00f79=7, 2\", \"26\u00f79=2, 8\"),\n    @(\"60\u00f73=20, 0\", \"57\u00f73=19, 0\"),\n    @(\"58\u00f76=9, 4\", \"73\u00f74=18, 1\"),\n    @(\"86\u00f73=28, 2\", \"46\u00f74=11, 2\"),\n    @(\"53\u00f79=5, 8\", \"12\u00f75=2, 2\"),\n    @(\"95\u00f73=31, 2\", \"95\u00f72=47, 1\"),\n    @(\"67\u00f72=33, 1\", \"61\u00f73=20, 1\"),\n    @(\"63\u00f79=7, 0\", \"17\u00f79=1, 8\"),\n    @(\"98\u00f79=10, 8\", \"86\u00f79=9, 5\"),\n    @(\"81\u00f79=9, 0\", \"24\u00f77=3, 3\"),\n    @(\"58\u00f73=19, 1\", \"75\u00f74=18, 3\"),\n    @(\"23\u00f79=2, 5\", \"56\u00f75=11, 1\"),\n    @(\"93\u00f74=23, 1\", \"23\u00f79=2, 5\"),\n    @(\"53\u00f72=26, 1\", \"93\u00f78=11, 5\"),\n    @(\"28\u00f75=5, 3\", \"10\u00f73=3, 1\"),\n    @(\"74\u00f76=12, 2\", \"27\u00f72=13, 1\"),\n    @(\"72\u00f78=9, 0\", \"49\u00f76=8, 1\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n    if (-not $found) {\n        throw \"Find/Replace failed for $oldText\"\n    }\n}\n\n"}
